$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '20.533.45'
$ws.Cells.Item(2, 5).Value = '  -0.10%  '
$ws.Cells.Item(3, 4).Value = '1.479.33'
$ws.Cells.Item(3, 5).Value = '  +0.56%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.009'
$ws.Cells.Item(4, 5).Value = '  -0.18%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.9767'
$ws.Cells.Item(5, 5).Value = '  -0.01%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '279.39'
$ws.Cells.Item(6, 5).Value = '  -0.85%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3659'
$ws.Cells.Item(7, 5).Value = '  -1.82%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3075'
$ws.Cells.Item(8, 5).Value = '  -4.10%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '39.91'
$ws.Cells.Item(9, 5).Value = '  -4.84%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.063'
$ws.Cells.Item(10, 5).Value = '  -0.67%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.06655'
$ws.Cells.Item(11, 5).Value = '  -1.08%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.004'
$ws.Cells.Item(12, 5).Value = '  -0.15%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.512'
$ws.Cells.Item(13, 5).Value = '  -2.41%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '18.15'
$ws.Cells.Item(14, 5).Value = '  -1.63%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.202'
$ws.Cells.Item(15, 5).Value = '  -1.72%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.9777'
$ws.Cells.Item(16, 5).Value = '  +0.62%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.00001031'
$ws.Cells.Item(17, 5).Value = '  -0.98%  '
$ws.Cells.Item(18, 4).Value = '1.480.74'
$ws.Cells.Item(18, 5).Value = '  +0.39%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.05935'
$ws.Cells.Item(19, 5).Value = '  +2.17%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '69.64'
$ws.Cells.Item(20, 5).Value = '  -4.66%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.482'
$ws.Cells.Item(21, 5).Value = '  -4.12%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '14.55'
$ws.Cells.Item(22, 5).Value = '  -2.67%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '11.04'
$ws.Cells.Item(23, 5).Value = '  -2.40%  '
$ws.Cells.Item(24, 5).Value = '  -2.60%  '
$ws.Cells.Item(25, 4).Value = '20.596.89'
$ws.Cells.Item(25, 5).Value = '  -0.08%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '141.14'
$ws.Cells.Item(26, 5).Value = '  +2.37%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.137'
$ws.Cells.Item(27, 5).Value = '  -8.82%  '
$ws.Cells.Item(28, 5).Value = '  -2.10%  '
$ws.Cells.Item(29, 4).Value = '1.639.08'
$ws.Cells.Item(29, 5).Value = '  -0.08%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '114.05'
$ws.Cells.Item(30, 5).Value = '  +0.13%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.996'
$ws.Cells.Item(31, 5).Value = '  +0.43%  '
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.995'
$ws.Cells.Item(32, 5).Value = '  -7.45%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.8138'
$ws.Cells.Item(33, 5).Value = '  -4.11%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.08022'
$ws.Cells.Item(34, 5).Value = '  +1.96%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.545'
$ws.Cells.Item(35, 5).Value = '  -5.50%  '
$ws.Cells.Item(36, 5).Value = '  +7.37%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.05843'
$ws.Cells.Item(37, 5).Value = '  -4.07%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.725'
$ws.Cells.Item(38, 5).Value = '  -4.85%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '7.795'
$ws.Cells.Item(39, 5).Value = '  +0.75%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.9771'
$ws.Cells.Item(40, 5).Value = '  -0.52%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.02047'
$ws.Cells.Item(41, 5).Value = '  -1.91%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '10.46'
$ws.Cells.Item(42, 5).Value = '  -3.31%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1886'
$ws.Cells.Item(43, 5).Value = '  -1.21%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.5300'
$ws.Cells.Item(44, 5).Value = '  -3.12%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.528'
$ws.Cells.Item(45, 5).Value = '  -2.03%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.14'
$ws.Cells.Item(46, 5).Value = '  -3.44%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '119.43'
$ws.Cells.Item(47, 5).Value = '  -2.07%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.5201'
$ws.Cells.Item(48, 5).Value = '  -3.29%  '
$ws.Cells.Item(49, 5).Value = '  -2.61%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.06465'
$ws.Cells.Item(50, 5).Value = '  +0.18%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.9913'
$ws.Cells.Item(51, 5).Value = '  -1.19%  '
